# Locator_Data.xlsx — "Add employee page in progress"
#
# The select_photo locator row (AddEmployeePage row 10) is split out of a
# single combined XPath expression into its three separate locator columns
# (cssselector / xpath / classname), matching the layout already used by
# every other row. The previously-active LoginPage selection moves to D4
# and AddEmployeePage becomes the active sheet/selected tab with its
# selection on E10.

$wb = $excel.ActiveWorkbook
$wsLogin = $wb.Worksheets.Item("LoginPage")
$wsAddEmp = $wb.Worksheets.Item("AddEmployeePage")

# select_photo locator: was a single xpath combining type+class;
# now split across cssselector (D), xpath (E), classname (F).
$wsAddEmp.Range("D10").Value = "input[type='file']"
$wsAddEmp.Range("E10").Value = "//input[@type='file']"
$wsAddEmp.Range("F10").Value = "oxd-file-input"

# LoginPage is no longer the active tab; its remembered selection moves to D4.
$wsLogin.Range("D4").Select()

# AddEmployeePage becomes the active / selected tab, with selection on E10.
$wsAddEmp.Activate()
$wsAddEmp.Range("E10").Select()
